# Generate Report for Handback
# Adds a new handback row (5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md) to the
# Overview / zh-cn / de-de sheets, expanding each sheet's table by one row.

$wb = $excel.ActiveWorkbook

$hyperColor = 15570276          # RGB(0x64,0x95,0xED) == FF6495ED used by the existing hyperlink cells
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Set-HandbackLink($ws, $cellRef, $address, $display) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $display) | Out-Null
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = $hyperColor
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
$wsOverview.Range("B4").Value = "e2e\5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-11-15 16:38:58"
$wsOverview.Range("G4").NumberFormat = $dateFormat

Set-HandbackLink $wsOverview "B4" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c083e564e5efb57416e4cbf8a4498fb32111abc1/e2e/5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md" `
    "e2e\5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.c083e564e5efb57416e4cbf8a4498fb32111abc1.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-11-15 16:38:44"
$wsZhCn.Range("H4").NumberFormat = $dateFormat
$wsZhCn.Range("I4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
$wsZhCn.Range("J4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.c083e564e5efb57416e4cbf8a4498fb32111abc1.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-11-15 16:39:39"
$wsZhCn.Range("K4").NumberFormat = $dateFormat
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

Set-HandbackLink $wsZhCn "A4" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c083e564e5efb57416e4cbf8a4498fb32111abc1/e2e/5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md" `
    "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
Set-HandbackLink $wsZhCn "I4" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c083e564e5efb57416e4cbf8a4498fb32111abc1/e2e/5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md" `
    "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.c083e564e5efb57416e4cbf8a4498fb32111abc1.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-11-15 16:38:58"
$wsDeDe.Range("H4").NumberFormat = $dateFormat
$wsDeDe.Range("I4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
$wsDeDe.Range("J4").Value = "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.c083e564e5efb57416e4cbf8a4498fb32111abc1.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-11-15 16:39:58"
$wsDeDe.Range("K4").NumberFormat = $dateFormat
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

Set-HandbackLink $wsDeDe "A4" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c083e564e5efb57416e4cbf8a4498fb32111abc1/e2e/5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md" `
    "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"
Set-HandbackLink $wsDeDe "I4" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c083e564e5efb57416e4cbf8a4498fb32111abc1/e2e/5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md" `
    "5a9ffbff-0fa4-45c0-b40d-510c6e44f825.md"

Write-Host "done"
